$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 12
$ws.Cells.Item(3, 6).Value = 343
$ws.Cells.Item(4, 6).Value = 1335
$ws.Cells.Item(5, 6).Value = 384
$ws.Cells.Item(7, 6).Value = 3930
$ws.Cells.Item(9, 6).Value = 788
$ws.Cells.Item(10, 6).Value = 2366
$ws.Cells.Item(11, 6).Value = 368
$ws.Cells.Item(12, 6).Value = 52
$ws.Cells.Item(14, 6).Value = 760
$ws.Cells.Item(15, 6).Value = 213
$ws.Cells.Item(16, 6).Value = 201
$ws.Cells.Item(17, 6).Value = 2991
$ws.Cells.Item(18, 6).Value = 326
$ws.Cells.Item(19, 6).Value = 237
$ws.Cells.Item(21, 6).Value = 353
$ws.Cells.Item(22, 6).Value = 246
$ws.Cells.Item(23, 6).Value = 51
$ws.Cells.Item(24, 6).Value = 282

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 57
$ws.Cells.Item(7, 6).Value = 135
$ws.Cells.Item(8, 6).Value = 105
$ws.Cells.Item(10, 6).Value = 102
$ws.Cells.Item(22, 6).Value = 80

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 834
$ws.Cells.Item(4, 6).Value = 2128
$ws.Cells.Item(5, 6).Value = 349
$ws.Cells.Item(6, 6).Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 834
$ws.Cells.Item(4, 6).Value = 2128
$ws.Cells.Item(5, 6).Value = 349
$ws.Cells.Item(6, 6).Value = 57
$ws.Cells.Item(7, 6).Value = 57
$ws.Cells.Item(8, 6).Value = 12
$ws.Cells.Item(10, 6).Value = 343
$ws.Cells.Item(11, 6).Value = 1335
$ws.Cells.Item(12, 6).Value = 384
$ws.Cells.Item(16, 6).Value = 19
$ws.Cells.Item(18, 6).Value = 3930
$ws.Cells.Item(19, 6).Value = 135
$ws.Cells.Item(21, 6).Value = 105
$ws.Cells.Item(23, 6).Value = 102
$ws.Cells.Item(24, 6).Value = 788
$ws.Cells.Item(25, 6).Value = 2366
$ws.Cells.Item(26, 6).Value = 368
$ws.Cells.Item(27, 6).Value = 52
$ws.Cells.Item(30, 6).Value = 760
$ws.Cells.Item(31, 6).Value = 213
$ws.Cells.Item(32, 6).Value = 201
$ws.Cells.Item(35, 6).Value = 326
$ws.Cells.Item(38, 6).Value = 237
$ws.Cells.Item(40, 6).Value = 353
$ws.Cells.Item(41, 6).Value = 246
$ws.Cells.Item(42, 6).Value = 51
$ws.Cells.Item(49, 6).Value = 80
$ws.Cells.Item(50, 6).Value = 282
